$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-02 Monday" "2024-09-03 Tuesday"

Replace-Text "48×70=3360" "58×70=4060"
Replace-Text "24×78=1872" "60×63=3780"
Replace-Text "58×60=3480" "94×45=4230"
Replace-Text "14×77=1078" "60×73=4380"
Replace-Text "83×40=3320" "71×83=5893"

Replace-Text "37×63=2331" "39×16=624"
Replace-Text "72×89=6408" "17×97=1649"
Replace-Text "89×52=4628" "97×48=4656"
Replace-Text "89×56=4984" "46×24=1104"
Replace-Text "75×75=5625" "50×44=2200"

Replace-Text "31×12=372" "66×40=2640"
Replace-Text "68×56=3808" "80×69=5520"
Replace-Text "24×38=912" "20×58=1160"
Replace-Text "68×55=3740" "76×97=7372"
Replace-Text "47×79=3713" "64×26=1664"

Replace-Text "58×29=1682" "94×49=4606"
Replace-Text "40×63=2520" "38×56=2128"
Replace-Text "11×45=495" "73×99=7227"
Replace-Text "94×75=7050" "80×61=4880"
Replace-Text "99×88=8712" "26×55=1430"

Replace-Text "86×82=7052" "57×50=2850"
Replace-Text "70×76=5320" "19×91=1729"
Replace-Text "31×57=1767" "93×65=6045"
Replace-Text "50×11=550" "62×19=1178"
Replace-Text "19×22=418" "62×98=6076"
